$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Group 1: rows 45-48 share the same note in column B, so column A was
# filled first (top to bottom) and then column B.
$ws.Range("A45").Value = "Anh sinh năm Mùi, Mệnh an tại Sửu"
$ws.Range("A46").Value = "Anh sinh năm Ngọ, Mệnh an tại Sửu"
$ws.Range("A47").Value = "Anh sinh năm Mùi, Mệnh an tại Tý"
$ws.Range("A48").Value = "Anh sinh năm Ngọ, Mệnh an tại Tý"

$ws.Range("B45").Value = "Cuộc đời vất vả,ít có điều xứng ý toại lòng."
$ws.Range("B46").Value = "Cuộc đời vất vả,ít có điều xứng ý toại lòng."
$ws.Range("B47").Value = "Cuộc đời vất vả,ít có điều xứng ý toại lòng."
$ws.Range("B48").Value = "Cuộc đời vất vả,ít có điều xứng ý toại lòng."

# Group 2: rows 49-50, entered one full row at a time.
$ws.Range("A49").Value = "Cung Mệnh của chị được an tại ví trí Tứ Mộ"
$ws.Range("B49").Value = "Chị là một người khôn ngoan, đảm đang."

$ws.Range("A50").Value = "Cung Mệnh của chị được an tại ví trí cung Dậu"
$ws.Range("B50").Value = "Chị đi ra ngoài được rất nhiều người  yêu mến và tôn trọng."

# Group 3: rows 51-52 share the same note in column B, column A filled
# first, then column B.
$ws.Range("A51").Value = "Cung Mệnh của chị được an tại ví trí cung Ngọ"
$ws.Range("A52").Value = "Cung Mệnh của chị được an tại ví trí cung Tý"

$ws.Range("B51").Value = "Cuộc đời chị an nhàn."
$ws.Range("B52").Value = "Cuộc đời chị an nhàn."

# Apply the highlighted style (yellow fill) used throughout column A to the
# new cells, matching the existing formatting pattern.
$ws.Range("A45:A52").Interior.Color = $ws.Range("A41").Interior.Color

# Update the view state to match the saved selection/scroll position.
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("E53").Select()
